$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D14").Value = 0.6447426901493167
$ws.Range("C15").Value = 0.2386249091493167
$ws.Range("D15").Value = 0.597740902
$ws.Range("B16").Value = -0.0107480648506833
$ws.Range("C16").Value = 0.042359665
